$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.697.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.95%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.503.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.44%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.77%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.40%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.494.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.14%  "

$ws.Range("E10").Value = "  -2.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.424"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.110.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.91%  "

$ws.Range("E15").Value = "  -1.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.715.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000173"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.497.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.548"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.39%  "

$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.33%  "

$ws.Range("E26").Value = "  +0.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000120"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.81%  "

$ws.Range("E29").Value = "  -2.05%  "

$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "24.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.02%  "

$ws.Range("E33").Value = "  -2.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.79%  "

$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "29.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "160.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.892"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.97%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.82%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.731.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0703"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.54%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0292"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "325.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.22%  "

$ws.Range("E51").Value = "  -3.77%  "
